$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - MegaScan
$ws.Range("A6").Value = "Roberto"
$ws.Range("B6").Value = "'0643"
$ws.Range("C6").Value = "MegaScan"
$ws.Range("D6").Value = "Sem comunicação de câmeras."
$ws.Range("G6").Value = "Pendente"

# Row 7 - Manoel Correira
$ws.Range("A7").Value = "Roberto"
$ws.Range("B7").Value = "'0756"
$ws.Range("C7").Value = "Manoel Correira"
$ws.Range("D7").Value = "Sem comunicação de câmeras, passar para o DDNS."
$ws.Range("G7").Value = "Pendente"

# Row 8 - Galpao Toyota (wraps to two lines)
$ws.Range("A8").Value = "Roberto"
$ws.Range("B8").Value = "'0803"
$ws.Range("C8").Value = "Galpao Toyota"
$ws.Range("D8").Value = "Disparos em falso, passar central pra internet e instalar/programar câmera pra gente monitorar."
$ws.Range("D8").WrapText = $true
$ws.Range("G8").Value = "Pendente"
$ws.Rows.Item(8).RowHeight = 30

# Row 9 - MedCenter
$ws.Range("A9").Value = "Roberto"
$ws.Range("B9").Value = "'0079"
$ws.Range("C9").Value = "MedCenter"
$ws.Range("D9").Value = "Sem comunicação de alarmes, local funciona via internet."
$ws.Range("G9").Value = "Pendente"

# Row 10 - Brapi
$ws.Range("A10").Value = "Roberto"
$ws.Range("B10").Value = "'0217"
$ws.Range("C10").Value = "Brapi"
$ws.Range("D10").Value = "Setores abertos, é AMT 8000. Pedro diretor pediu pra ver sobre. "
$ws.Range("G10").Value = "Pendente"

# Restore the active selection to C8, matching the author's last cursor position
$ws.Range("C8").Select()
